$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "Table_AS_2023-11-13"
$wb.Worksheets.Item(2).Name = "Table_CT_2023-11-13"
$wb.Worksheets.Item(3).Name = "Table_AS-CT_2023-11-13"

# --- Add a fifth "reporting" column to the two sheets that contain mapping rows ---
foreach ($idx in 1, 3) {
    $ws = $wb.Worksheets.Item($idx)

    # Find the last used row in column A
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    # Header: clone the style of the existing header cell (D1) onto E1, then set its text
    $ws.Range("D1").Copy() | Out-Null
    $ws.Range("E1").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("E1").Value = "reporting"

    # Data rows: "<n>. '<olabel>' --> '<slabel>'"
    for ($r = 2; $r -le $lastRow; $r++) {
        $oLabel = $ws.Cells.Item($r, 1).Value()
        $sLabel = $ws.Cells.Item($r, 3).Value()
        $n = $r - 1
        $ws.Cells.Item($r, 5).Value = "$n. '$oLabel' --> '$sLabel'"
    }
}
